$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3045.125
$ws.Range("I2").Value = 317.5
$ws.Range("J2").Value = 5772.75
$ws.Range("K2").Value = 317.5
$ws.Range("L2").Value = 5772.75
$ws.Range("M2").Value = -204.5
$ws.Range("N2").Value = -5998.75

$ws.Range("H33").Value = 307.0909
$ws.Range("I33").Value = 322.5
$ws.Range("K33").Value = 322.5
$ws.Range("M33").Value = -93.5

$ws.Range("H51").Value = 5552
$ws.Range("J51").Value = 5552
$ws.Range("L51").Value = 5552
$ws.Range("N51").Value = -6520

$ws.Range("H64").Value = 4829.6665
$ws.Range("J64").Value = 5500
$ws.Range("L64").Value = 5500
$ws.Range("N64").Value = -5996

$ws.Range("H67").Value = 4829.6665
$ws.Range("J67").Value = 5500
$ws.Range("L67").Value = 5500
$ws.Range("N67").Value = -7216

$ws.Range("H74").Value = 27642.111
$ws.Range("I74").Value = 30059.834
$ws.Range("K74").Value = 30059.834
$ws.Range("M74").Value = -29123.834

$ws.Range("H77").Value = 27642.111
$ws.Range("I77").Value = 30059.834
$ws.Range("K77").Value = 150299.17
$ws.Range("M77").Value = -145619.17

$ws.Range("H86").Value = 949.3333
$ws.Range("I86").Value = 659.8
$ws.Range("J86").Value = 1311.25
$ws.Range("K86").Value = 659.8
$ws.Range("L86").Value = 1311.25
$ws.Range("M86").Value = 463.2
$ws.Range("N86").Value = -3557.25

$ws.Range("H89").Value = 949.3333
$ws.Range("I89").Value = 659.8
$ws.Range("J89").Value = 1311.25
$ws.Range("K89").Value = 3299
$ws.Range("L89").Value = 6556.25
$ws.Range("M89").Value = 2317
$ws.Range("N89").Value = -17788.25

$ws.Range("H99").Value = 1508
$ws.Range("I99").Value = 525
$ws.Range("J99").Value = 1999.5
$ws.Range("K99").Value = 1575
$ws.Range("L99").Value = 5998.5
$ws.Range("M99").Value = -77
$ws.Range("N99").Value = -8994.5

$ws.Range("I101").Value = 14290192
$ws.Range("K101").Value = 42870576
$ws.Range("M101").Value = -42868954

$ws.Range("H116").Value = 3824.4333
$ws.Range("I116").Value = 2949.9285
$ws.Range("K116").Value = 2949.9285
$ws.Range("M116").Value = 492.0715

$ws.Range("H118").Value = 1329.8334
$ws.Range("I118").Value = 995.8
$ws.Range("K118").Value = 2987.4
$ws.Range("M118").Value = -1330.4

$ws.Range("H132").Value = 3021.4443
$ws.Range("I132").Value = 3070.4285
$ws.Range("J132").Value = 2850
$ws.Range("K132").Value = 9211.2855
$ws.Range("L132").Value = 8550
$ws.Range("M132").Value = -6681.2855
$ws.Range("N132").Value = -13610

$ws.Range("H135").Value = 1947.4694
$ws.Range("I135").Value = 2087.15
$ws.Range("K135").Value = 18784.35
$ws.Range("M135").Value = -16249.35

$ws.Range("H137").Value = 1166.2778
$ws.Range("I137").Value = 954.1667
$ws.Range("J137").Value = 1590.5
$ws.Range("K137").Value = 2862.5001
$ws.Range("L137").Value = 4771.5
$ws.Range("M137").Value = -312.5001000000002
$ws.Range("N137").Value = -9871.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3519.7144
$ws.Range("I45").Value = 3519.7144
$ws.Range("K45").Value = 3519.7144
$ws.Range("M45").Value = -3142.7144

$ws.Range("H46").Value = 9711
$ws.Range("J46").Value = 9698
$ws.Range("L46").Value = 9698
$ws.Range("N46").Value = -10336

$ws.Range("H61").Value = 1959.9333
$ws.Range("J61").Value = 1625.3334
$ws.Range("L61").Value = 1625.3334
$ws.Range("N61").Value = -2049.3334

$ws.Range("H97").Value = 868.25
$ws.Range("I97").Value = 740
$ws.Range("K97").Value = 740
$ws.Range("M97").Value = -244

$ws.Range("H136").Value = 1959.9333
$ws.Range("J136").Value = 1625.3334
$ws.Range("L136").Value = 4876.0002
$ws.Range("N136").Value = -9976.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 550
$ws.Range("I5").Value = 550
$ws.Range("K5").Value = 550
$ws.Range("M5").Value = -437

$ws.Range("H20").Value = 9001.75
$ws.Range("I20").Value = 8669
$ws.Range("K20").Value = 8669
$ws.Range("M20").Value = -8422

$ws.Range("H119").Value = 761
$ws.Range("J119").Value = 761
$ws.Range("L119").Value = 761
$ws.Range("N119").Value = -10437

$ws.Range("H134").Value = 2429.5
$ws.Range("I134").Value = 2499.4443
$ws.Range("J134").Value = 1800
$ws.Range("K134").Value = 7498.3329
$ws.Range("L134").Value = 5400
$ws.Range("M134").Value = -4963.3329
$ws.Range("N134").Value = -10470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2224.6667
$ws.Range("I31").Value = 1673.3
$ws.Range("J31").Value = 2549
$ws.Range("K31").Value = 1673.3
$ws.Range("L31").Value = 2549
$ws.Range("M31").Value = -1378.3
$ws.Range("N31").Value = -3139

$ws.Range("H34").Value = 2224.6667
$ws.Range("I34").Value = 1673.3
$ws.Range("J34").Value = 2549
$ws.Range("K34").Value = 1673.3
$ws.Range("L34").Value = 2549
$ws.Range("M34").Value = -1471.3
$ws.Range("N34").Value = -2953

$ws.Range("H132").Value = 1667.9166
$ws.Range("I132").Value = 1667.9166
$ws.Range("K132").Value = 5003.7498
$ws.Range("M132").Value = -2473.7498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3067.4707
$ws.Range("J68").Value = 3276.5334
$ws.Range("L68").Value = 9829.600199999999
$ws.Range("N68").Value = -11451.6002

$ws.Range("H71").Value = 3067.4707
$ws.Range("J71").Value = 3276.5334
$ws.Range("L71").Value = 29488.8006
$ws.Range("N71").Value = -37600.8006

$ws.Range("H87").Value = 13985
$ws.Range("I87").Value = 13985
$ws.Range("K87").Value = 41955
$ws.Range("M87").Value = -40707

$ws.Range("H90").Value = 13985
$ws.Range("I90").Value = 13985
$ws.Range("K90").Value = 125865
$ws.Range("M90").Value = -119625

$ws.Range("H107").Value = 1870.6364
$ws.Range("J107").Value = 2117.2
$ws.Range("L107").Value = 6351.599999999999
$ws.Range("N107").Value = -10191.6

$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = ""

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = ""

$ws.Range("H112").Value = 58999
$ws.Range("J112").Value = 58999
$ws.Range("L112").Value = 58999
$ws.Range("N112").Value = -61215

$ws.Range("H113").Value = 874.5714
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""

$ws.Range("H126").Value = 10010
$ws.Range("I126").Value = 10000
$ws.Range("J126").Value = 10013.333
$ws.Range("K126").Value = 30000
$ws.Range("L126").Value = 30039.999
$ws.Range("M126").Value = -27530
$ws.Range("N126").Value = -34979.999

$ws.Range("H132").Value = 5666.6665
$ws.Range("I132").Value = 5666.6665
$ws.Range("K132").Value = 16999.9995
$ws.Range("M132").Value = -14469.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3120.3635
$ws.Range("I46").Value = 1642.3334
$ws.Range("J46").Value = 4894
$ws.Range("K46").Value = 1642.3334
$ws.Range("L46").Value = 4894
$ws.Range("M46").Value = -1454.3334
$ws.Range("N46").Value = -5270

$ws.Range("H93").Value = 643.8182
$ws.Range("I93").Value = 659.2
$ws.Range("K93").Value = 659.2
$ws.Range("M93").Value = 588.8

$ws.Range("H106").Value = 47067.5
$ws.Range("J106").Value = 47067.5
$ws.Range("L106").Value = 47067.5
$ws.Range("N106").Value = -49591.5

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""

$ws.Range("H100").Value = 20001674
$ws.Range("I100").Value = 25001924
$ws.Range("K100").Value = 50003848
$ws.Range("M100").Value = -50003307

$ws.Range("H105").Value = 9000
$ws.Range("J105").Value = 9000
$ws.Range("L105").Value = 9000
$ws.Range("N105").Value = -15988

$ws.Range("H132").Value = 550
$ws.Range("I132").Value = 550
$ws.Range("K132").Value = 1650
$ws.Range("M132").Value = 880

$ws.Range("H136").Value = 2416.0557
$ws.Range("I136").Value = 2300.3333
$ws.Range("J136").Value = 2647.5
$ws.Range("K136").Value = 6900.999899999999
$ws.Range("L136").Value = 7942.5
$ws.Range("M136").Value = -4350.999899999999
$ws.Range("N136").Value = -13042.5
